$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update stats for 2025-09 (row 22)
$ws.Range("B22").Value = 6288
$ws.Range("C22").Value = 995
$ws.Range("D22").Value = 5838658
$ws.Range("E22").Value = 928.5397582697201
$ws.Range("F22").Value = 8.245825443277678
$ws.Range("G22").Value = 4.079497907949792
$ws.Range("H22").Value = 26.97261632408077
